# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-23 14:12:35
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Row 2: re-order the "Recorded By" list for the ANATOMY session 1
$ws.Range("G2").Value = "System, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# Row 3: re-order the "Recorded By" list for the ANATOMY session 2
$ws.Range("G3").Value = "System, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# Row 4: re-order the "Recorded By" list for the ANATOMY session 3
$ws.Range("G4").Value = "hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# Row 5: ANATOMY session 4 - a new recorder was added and attendance count increased
$ws.Range("G5").Value = "Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("H5").Value = "54/251"

# Row 7: re-order the "Recorded By" list for the BIOCHEMISTRY LAB/CBL session 1
$ws.Range("G7").Value = "Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"

# Row 9: re-order the "Recorded By" list for the HISTOLOGY session 1
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

# Updated "Average Attendance %" metric (recomputed after the row 5 change above).
# Force text so Excel doesn't reinterpret the literal percent string as a number.
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "23.5%"

# Same "Average Attendance %" figure repeated in the per-subgroup summary table
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "23.5%"
